# Scene.xlsx optimization:
#  - Row 2 now describes the "CloneScene" entry (ID=3, SceneName="clone")
#    instead of the old "Login"/ID=0 entry.
#  - Rows 3 and 4 switch their SceneName from "Stage001" to "newscene".
#  - Active selection moves to H8.
#
# NOTE: set B2 (the ID column, stored as text) before A2 so that the
# shared-string table gets new entries appended in the same order Excel
# produced them in the source workbook ("3" before the CloneScene path).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "3"
$ws.Range("A2").Value = "../../NFDataCfg/Ini/NFZoneServer/Scene/CloneScene/"
$ws.Range("F2").Value = "clone"
$ws.Range("F3").Value = "newscene"
$ws.Range("F4").Value = "newscene"

[void]$ws.Range("H8").Select()
